$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 20): date in column A, count in column B
# Copy the format from the cell above (A19) so the new date cell uses the
# same existing date style instead of creating a new number format.
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A20").Value = 45901
$ws.Range("B20").Value = 6

# Update the active selection to match the authored state
$ws.Range("D19").Select()
